$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.918.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5063"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06451"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.281"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.875.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.646.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5622"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7729"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.964.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.370"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.953"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.128"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("E25").Value = "  -5.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1241"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.819"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04977"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.307"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.242"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.578"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.389"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9074"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5580"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.566"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.131.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01569"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.527"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8030"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "98.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.784.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4296"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.794"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05050"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9994"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.65%  "
